# Re-executed Ledger and View Register page
# - Re-encrypted "password" crypt value (row 21)
# - Updated Type.EnterValue test data: "Purchase" -> "Sale" (row 73)
# - Appended new "View Register" (vr.*) field/xpath rows 90-108
# - Select the last entered cell to mirror the recorded cursor position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Re-encrypted password value
$ws.Range("B21").Value = 'crypt:e017e92383a6feabc0ed331d068d7fea4f207739bd83bcf1'
$ws.Range("B21").Font.Bold = $true

# 2) Test data value changed from Purchase to Sale
$ws.Range("B73").Value = 'Sale'

# 3) New View Register rows (vr.* labels + their xpaths)
$ws.Range("A90").Value = 'vr.No'
$ws.Range("B90").Value = '//main[@class=''mb-5'']/section/div[2]/div/table/thead/tr/th[1]/div'
$ws.Range("A91").Value = 'vr.Inv.No'
$ws.Range("B91").Value = '//main[@class=''mb-5'']/section/div[2]/div/table/thead/tr/th[2]/div'
$ws.Range("A92").Value = 'vr.Inv.Date'
$ws.Range("B92").Value = '//main[@class=''mb-5'']/section/div[2]/div/table/thead/tr/th[3]/div'
$ws.Range("A93").Value = 'vr.Customer'
$ws.Range("B93").Value = '//main[@class=''mb-5'']/section/div[2]/div/table/thead/tr/th[4]/div'
$ws.Range("A94").Value = 'vr.GSTIN'
$ws.Range("B94").Value = '//main[@class=''mb-5'']/section/div[2]/div/table/thead/tr/th[5]/div'
$ws.Range("A95").Value = 'vr.TIN'
$ws.Range("B95").Value = '//main[@class=''mb-5'']/section/div[2]/div/table/thead/tr/th[6]/div'
$ws.Range("A96").Value = 'vr.TaxFree'
$ws.Range("B96").Value = '//main[@class=''mb-5'']/section/div[2]/div/table/thead/tr/th[7]/div'
$ws.Range("A97").Value = 'vr.Taxable'
$ws.Range("B97").Value = '//main[@class=''mb-5'']/section/div[2]/div/table/thead/tr/th[8]/div'
$ws.Range("A98").Value = 'vr.GAmount'
$ws.Range("B98").Value = '//main[@class=''mb-5'']/section/div[2]/div/table/thead/tr/th[9]/div'
$ws.Range("A99").Value = 'vr.table.row.xpath'
$ws.Range("B99").Value = '//main[@class=''mb-5'']/section/div[2]/div/table/tbody/tr'
$ws.Range("A100").Value = 'vr.No.values'
$ws.Range("B100").Value = '/td[1]'
$ws.Range("A101").Value = 'vr.Inv.No.values'
$ws.Range("B101").Value = '/td[2]'
$ws.Range("A102").Value = 'vr.Inv.Date.values'
$ws.Range("B102").Value = '/td[3]'
$ws.Range("A103").Value = 'vr.Customer.values'
$ws.Range("B103").Value = '/td[4]'
$ws.Range("A104").Value = 'vr.GSTIN.values'
$ws.Range("B104").Value = '/td[5]'
$ws.Range("A105").Value = 'vr.TIN.values'
$ws.Range("B105").Value = '/td[6]'
$ws.Range("A106").Value = 'vr.TaxFree.values'
$ws.Range("B106").Value = '/td[7]'
$ws.Range("A107").Value = 'vr.Taxable.values'
$ws.Range("B107").Value = '/td[8]'
$ws.Range("A108").Value = 'vr.GAmount.values'
$ws.Range("B108").Value = '/td[9]'

# 4) Mirror the recorded selection after appending the new rows
$ws.Range("A108").Select()
